$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 3 (shifts existing rows 3..22 down to 4..23),
# producing a new data point for quarter 2020-04-01 that the averaging
# had previously skipped (selection-scope bug being fixed here).
$ws.Rows(3).Insert()

# Copy the label cell formatting (bold, centered, bordered) from the row
# below (the old row 3, now row 4) onto the freshly inserted label cell.
$ws.Range("A4").Copy()
$ws.Range("A3").PasteSpecial(-4122)

# New row label: the quarter that was missing from the averaged series.
$ws.Range("A3").Value = "2020-04-01 00:00:00_diff"

# New averaged-error values for the restored quarter.
$ws.Range("B3").Value = 7.469150330857293
$ws.Range("C3").Value = -10.70211146928018
$ws.Range("D3").Value = -2.562439922301026
$ws.Range("E3").Value = -0.4148694505016339
$ws.Range("F3").Value = -3.854899437024964
$ws.Range("G3").Value = -3.916073615439165
$ws.Range("H3").Value = -1.87461767828291
